# Add a new trade row (row 6) to the CELG named-trade worksheet, matching
# the style/format of the preceding row (row 5), and widen column E so the
# new BuyPrice value fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats/styles) of row 5 down into row 6 first,
# so that the date cell (A6) and the boolean-styled cell (G6) pick up the
# same cell style (s="1") that rows 3-5 already use.
$ws.Range("A5:I5").Copy() | Out-Null
$ws.Range("A6:I6").PasteSpecial(-4122) | Out-Null

# Now populate the new row's values.
$ws.Range("A6").Value = 42650.371504629627
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 9906.81
$ws.Range("D6").Value = 9949.09
$ws.Range("E6").Value = 104.839996
$ws.Range("F6").Value = 103.95
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = -0.85
$ws.Range("I6").Value = $false

# Widen column E (BuyPrice) so the new, wider value fits -- matches the
# bestFit recalculation Excel performs automatically for this column.
$ws.Columns.Item(5).ColumnWidth = 10.0
